$wb = $excel.ActiveWorkbook

# --- Add the two new worksheets at the end of the workbook ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ahaGlucose = $wb.Worksheets.Add($null, $lastSheet)
$ahaGlucose.Name = "aha glucose"
$ahaPatient = $wb.Worksheets.Add($null, $ahaGlucose)
$ahaPatient.Name = "aha patient"

# --- Populate "aha glucose" ---
$ahaGlucose.Range("A1").Value = "variable"
$ahaGlucose.Range("B1").Value = "new_var"
$ahaGlucose.Range("A2").Value = "Record ID   Group: [screeningrandomiza_arm_1][bllinded_group]"
$ahaGlucose.Range("B2").Value = "record_id"
$ahaGlucose.Range("A3").Value = "Event Name"
$ahaGlucose.Range("B3").Value = "event_name"
$ahaGlucose.Range("A4").Value = "Length of surgery"
$ahaGlucose.Range("B4").Value = "duration_surgery"
$ahaGlucose.Range("A5").Value = "Surgery start time"
$ahaGlucose.Range("B5").Value = "surgery_start_time"
$ahaGlucose.Range("A6").Value = "Surgery End Time"
$ahaGlucose.Range("B6").Value = "surgery_end_time"
$ahaGlucose.Range("A7").Value = "Date of CABG (BG monitoring)"
$ahaGlucose.Range("B7").Value = "date_cabg"
$ahaGlucose.Range("A8").Value = "OR Time 1:"
$ahaGlucose.Range("B8").Value = "time_or1"
$ahaGlucose.Range("A9").Value = "OR BG...8"
$ahaGlucose.Range("B9").Value = "glucose_or1"
$ahaGlucose.Range("A10").Value = "OR CGM value -Time 1"
$ahaGlucose.Range("A11").Value = "OR CGM Value1"
$ahaGlucose.Range("A12").Value = "OR Time 2:"
$ahaGlucose.Range("B12").Value = "time_or2"
$ahaGlucose.Range("A13").Value = "OR BG...12"
$ahaGlucose.Range("B13").Value = "glucose_or2"
$ahaGlucose.Range("A14").Value = "OR CGM -Time2"
$ahaGlucose.Range("A15").Value = "OR CGM Value2"
$ahaGlucose.Range("A16").Value = "OR Time 3:"
$ahaGlucose.Range("B16").Value = "time_or3"
$ahaGlucose.Range("A17").Value = "OR BG...16"
$ahaGlucose.Range("B17").Value = "glucose_or3"
$ahaGlucose.Range("A18").Value = "OR CGM -Time3"
$ahaGlucose.Range("A19").Value = "OR CGM Value3"
$ahaGlucose.Range("A20").Value = "OR Time 4:"
$ahaGlucose.Range("B20").Value = "time_or4"
$ahaGlucose.Range("A21").Value = "OR BG...20"
$ahaGlucose.Range("B21").Value = "glucose_or4"
$ahaGlucose.Range("A22").Value = "OR CGM -Time 4"
$ahaGlucose.Range("A23").Value = "OR CGM Value 4"
$ahaGlucose.Range("A24").Value = "OR Time 5:"
$ahaGlucose.Range("B24").Value = "time_or5"
$ahaGlucose.Range("A25").Value = "OR BG...24"
$ahaGlucose.Range("B25").Value = "glucose_or5"
$ahaGlucose.Range("A26").Value = "OR CGM -Time 5"
$ahaGlucose.Range("A27").Value = "OR CGM Value 5"
$ahaGlucose.Range("A28").Value = "OR Time 6:"
$ahaGlucose.Range("B28").Value = "time_or6"
$ahaGlucose.Range("A29").Value = "OR BG...28"
$ahaGlucose.Range("B29").Value = "glucose_or6"
$ahaGlucose.Range("A30").Value = "OR CGM -Time 6"
$ahaGlucose.Range("A31").Value = "OR CGM Value 6"
$ahaGlucose.Range("A32").Value = "OR Time 7:"
$ahaGlucose.Range("B32").Value = "time_or7"
$ahaGlucose.Range("A33").Value = "OR BG...32"
$ahaGlucose.Range("B33").Value = "glucose_or7"
$ahaGlucose.Range("A34").Value = "OR CGM -Time 7"
$ahaGlucose.Range("A35").Value = "OR CGM Value 7"
$ahaGlucose.Range("A36").Value = "OR Time 8:"
$ahaGlucose.Range("B36").Value = "time_or8"
$ahaGlucose.Range("A37").Value = "OR CGM -Time 8"
$ahaGlucose.Range("A38").Value = "OR CGM Value 8"
$ahaGlucose.Range("A39").Value = "OR BG...38"
$ahaGlucose.Range("B39").Value = "glucose_or8"
$ahaGlucose.Range("A40").Value = "OR Time 9:"
$ahaGlucose.Range("B40").Value = "time_or9"
$ahaGlucose.Range("A41").Value = "OR BG...40"
$ahaGlucose.Range("B41").Value = "glucose_or9"
$ahaGlucose.Range("A42").Value = "OR CGM -Time 9"
$ahaGlucose.Range("A43").Value = "OR CGM Value 9"
$ahaGlucose.Range("A44").Value = "OR Time 10:"
$ahaGlucose.Range("B44").Value = "time_or10"
$ahaGlucose.Range("A45").Value = "OR BG...44"
$ahaGlucose.Range("B45").Value = "glucose_or10"
$ahaGlucose.Range("A46").Value = "OR CGM -Time 10"
$ahaGlucose.Range("A47").Value = "OR CGM Value 10"
$ahaGlucose.Range("A48").Value = "OR Time 11:"
$ahaGlucose.Range("B48").Value = "time_or11"
$ahaGlucose.Range("A49").Value = "OR BG...48"
$ahaGlucose.Range("B49").Value = "glucose_or11"
$ahaGlucose.Range("A50").Value = "OR CGM -Time 11"
$ahaGlucose.Range("A51").Value = "OR CGM Value 11"
$ahaGlucose.Range("A52").Value = "OR Time 12:"
$ahaGlucose.Range("B52").Value = "time_or12"
$ahaGlucose.Range("A53").Value = "OR BG...52"
$ahaGlucose.Range("B53").Value = "glucose_or12"
$ahaGlucose.Range("A54").Value = "OR CGM -Time 12"
$ahaGlucose.Range("A55").Value = "OR CGM Value 12"
$ahaGlucose.Range("A56").Value = "OR Time 13:"
$ahaGlucose.Range("B56").Value = "time_or13"
$ahaGlucose.Range("A57").Value = "OR BG...56"
$ahaGlucose.Range("B57").Value = "glucose_or13"
$ahaGlucose.Range("A58").Value = "OR CGM -Time 13"
$ahaGlucose.Range("A59").Value = "OR CGM Value 13"
$ahaGlucose.Range("A60").Value = "OR Time 14:"
$ahaGlucose.Range("B60").Value = "time_or14"
$ahaGlucose.Range("A61").Value = "OR BG...60"
$ahaGlucose.Range("B61").Value = "glucose_or14"
$ahaGlucose.Range("A62").Value = "OR CGM -Time 14"
$ahaGlucose.Range("A63").Value = "OR CGM Value 14"
$ahaGlucose.Range("A64").Value = "OR Time 15:"
$ahaGlucose.Range("B64").Value = "time_or15"
$ahaGlucose.Range("A65").Value = "OR BG...64"
$ahaGlucose.Range("B65").Value = "glucose_or15"
$ahaGlucose.Range("A66").Value = "OR Time 16:"
$ahaGlucose.Range("B66").Value = "time_or16"
$ahaGlucose.Range("A67").Value = "OR BG...66"
$ahaGlucose.Range("B67").Value = "glucose_or16"
$ahaGlucose.Range("A68").Value = "OR Time 17:"
$ahaGlucose.Range("B68").Value = "time_or17"
$ahaGlucose.Range("A69").Value = "OR BG...68"
$ahaGlucose.Range("B69").Value = "glucose_or17"
$ahaGlucose.Range("A70").Value = "OR Time 18:"
$ahaGlucose.Range("B70").Value = "time_or18"
$ahaGlucose.Range("A71").Value = "OR BG...70"
$ahaGlucose.Range("B71").Value = "glucose_or18"
$ahaGlucose.Range("A72").Value = "OR Time 19:"
$ahaGlucose.Range("B72").Value = "time_or19"
$ahaGlucose.Range("A73").Value = "OR BG...72"
$ahaGlucose.Range("B73").Value = "glucose_or19"
$ahaGlucose.Range("A74").Value = "OR Time 20:"
$ahaGlucose.Range("B74").Value = "time_or20"
$ahaGlucose.Range("A75").Value = "OR BG...74"
$ahaGlucose.Range("B75").Value = "glucose_or20"
$ahaGlucose.Range("A76").Value = "Number of BG values"
$ahaGlucose.Range("A77").Value = "Average BG During Surgery (OR)"
$ahaGlucose.Range("A78").Value = "HyperglycemiaNumber of BG >=140mg/dL"
$ahaGlucose.Range("A79").Value = "Hyperglycemia BG >180mg/dL"
$ahaGlucose.Range("A80").Value = "Hypoglycemia BG < =70mg/dL and >40mg/dL"
$ahaGlucose.Range("A81").Value = "Hypoglycemia BG < =40mg/dL"
$ahaGlucose.Range("A82").Value = "Complete?"
$ahaGlucose.Range("A83").Value = "BG post-op"
$ahaGlucose.Range("B83").Value = "glucose_postop"
$ahaGlucose.Range("A84").Value = "Time 1:"
$ahaGlucose.Range("B84").Value = "time_icu1"
$ahaGlucose.Range("A85").Value = "BG...84"
$ahaGlucose.Range("B85").Value = "glucose_icu1"
$ahaGlucose.Range("A86").Value = "Time 2:"
$ahaGlucose.Range("B86").Value = "time_icu2"
$ahaGlucose.Range("A87").Value = "BG...86"
$ahaGlucose.Range("B87").Value = "glucose_icu2"
$ahaGlucose.Range("A88").Value = "Time 3:"
$ahaGlucose.Range("B88").Value = "time_icu3"
$ahaGlucose.Range("A89").Value = "BG...88"
$ahaGlucose.Range("B89").Value = "glucose_icu3"
$ahaGlucose.Range("A90").Value = "Time 4:"
$ahaGlucose.Range("B90").Value = "time_icu4"
$ahaGlucose.Range("A91").Value = "BG...90"
$ahaGlucose.Range("B91").Value = "glucose_icu4"
$ahaGlucose.Range("A92").Value = "Time 5:"
$ahaGlucose.Range("B92").Value = "time_icu5"
$ahaGlucose.Range("A93").Value = "BG...92"
$ahaGlucose.Range("B93").Value = "glucose_icu5"
$ahaGlucose.Range("A94").Value = "Time 6:"
$ahaGlucose.Range("B94").Value = "time_icu6"
$ahaGlucose.Range("A95").Value = "BG...94"
$ahaGlucose.Range("B95").Value = "glucose_icu6"
$ahaGlucose.Range("A96").Value = "Time 7:"
$ahaGlucose.Range("B96").Value = "time_icu7"
$ahaGlucose.Range("A97").Value = "BG...96"
$ahaGlucose.Range("B97").Value = "glucose_icu7"
$ahaGlucose.Range("A98").Value = "Time 8:"
$ahaGlucose.Range("B98").Value = "time_icu8"
$ahaGlucose.Range("A99").Value = "BG...98"
$ahaGlucose.Range("B99").Value = "glucose_icu8"
$ahaGlucose.Range("A100").Value = "Time 9:"
$ahaGlucose.Range("B100").Value = "time_icu9"
$ahaGlucose.Range("A101").Value = "BG...100"
$ahaGlucose.Range("B101").Value = "glucose_icu9"
$ahaGlucose.Range("A102").Value = "Time 10:"
$ahaGlucose.Range("B102").Value = "time_icu10"
$ahaGlucose.Range("A103").Value = "BG...102"
$ahaGlucose.Range("B103").Value = "glucose_icu10"
$ahaGlucose.Range("A104").Value = "Time 11:"
$ahaGlucose.Range("B104").Value = "time_icu11"
$ahaGlucose.Range("A105").Value = "BG...104"
$ahaGlucose.Range("B105").Value = "glucose_icu11"
$ahaGlucose.Range("A106").Value = "Time 12:"
$ahaGlucose.Range("B106").Value = "time_icu12"
$ahaGlucose.Range("A107").Value = "BG...106"
$ahaGlucose.Range("B107").Value = "glucose_icu12"
$ahaGlucose.Range("A108").Value = "Time 13:"
$ahaGlucose.Range("B108").Value = "time_icu13"
$ahaGlucose.Range("A109").Value = "BG...108"
$ahaGlucose.Range("B109").Value = "glucose_icu13"
$ahaGlucose.Range("A110").Value = "Time 14:"
$ahaGlucose.Range("B110").Value = "time_icu14"
$ahaGlucose.Range("A111").Value = "BG...110"
$ahaGlucose.Range("B111").Value = "glucose_icu14"
$ahaGlucose.Range("A112").Value = "Time 15:"
$ahaGlucose.Range("B112").Value = "time_icu15"
$ahaGlucose.Range("A113").Value = "BG...112"
$ahaGlucose.Range("B113").Value = "glucose_icu15"
$ahaGlucose.Range("A114").Value = "Time 16:"
$ahaGlucose.Range("B114").Value = "time_icu16"
$ahaGlucose.Range("A115").Value = "BG...114"
$ahaGlucose.Range("B115").Value = "glucose_icu16"
$ahaGlucose.Range("A116").Value = "Time 17:"
$ahaGlucose.Range("B116").Value = "time_icu17"
$ahaGlucose.Range("A117").Value = "BG...116"
$ahaGlucose.Range("B117").Value = "glucose_icu17"
$ahaGlucose.Range("A118").Value = "Time 18:"
$ahaGlucose.Range("B118").Value = "time_icu18"
$ahaGlucose.Range("A119").Value = "BG...118"
$ahaGlucose.Range("B119").Value = "glucose_icu18"
$ahaGlucose.Range("A120").Value = "Time 19:"
$ahaGlucose.Range("B120").Value = "time_icu19"
$ahaGlucose.Range("A121").Value = "BG...120"
$ahaGlucose.Range("B121").Value = "glucose_icu19"
$ahaGlucose.Range("A122").Value = "Time 20:"
$ahaGlucose.Range("B122").Value = "time_icu20"
$ahaGlucose.Range("A123").Value = "BG...122"
$ahaGlucose.Range("B123").Value = "glucose_icu20"
$ahaGlucose.Range("A124").Value = "Time 21:"
$ahaGlucose.Range("B124").Value = "time_icu21"
$ahaGlucose.Range("A125").Value = "BG...124"
$ahaGlucose.Range("B125").Value = "glucose_icu21"
$ahaGlucose.Range("A126").Value = "Time 22:"
$ahaGlucose.Range("B126").Value = "time_icu22"
$ahaGlucose.Range("A127").Value = "BG...126"
$ahaGlucose.Range("B127").Value = "glucose_icu22"
$ahaGlucose.Range("A128").Value = "Time 23:"
$ahaGlucose.Range("B128").Value = "time_icu23"
$ahaGlucose.Range("A129").Value = "BG...128"
$ahaGlucose.Range("B129").Value = "glucose_icu23"
$ahaGlucose.Range("A130").Value = "Time 24:"
$ahaGlucose.Range("B130").Value = "time_icu24"
$ahaGlucose.Range("A131").Value = "BG...130"
$ahaGlucose.Range("B131").Value = "glucose_icu24"
$ahaGlucose.Range("A132").Value = "Time 25:"
$ahaGlucose.Range("B132").Value = "time_icu25"
$ahaGlucose.Range("A133").Value = "BG...132"
$ahaGlucose.Range("B133").Value = "glucose_icu25"
$ahaGlucose.Range("A134").Value = "Time 26:"
$ahaGlucose.Range("B134").Value = "time_icu26"
$ahaGlucose.Range("A135").Value = "BG...134"
$ahaGlucose.Range("B135").Value = "glucose_icu26"
$ahaGlucose.Range("A136").Value = "Time 27:"
$ahaGlucose.Range("B136").Value = "time_icu27"
$ahaGlucose.Range("A137").Value = "BG...136"
$ahaGlucose.Range("B137").Value = "glucose_icu27"
$ahaGlucose.Range("A138").Value = "Time 28:"
$ahaGlucose.Range("B138").Value = "time_icu28"
$ahaGlucose.Range("A139").Value = "BG...138"
$ahaGlucose.Range("B139").Value = "glucose_icu28"
$ahaGlucose.Range("A140").Value = "Time 29:"
$ahaGlucose.Range("B140").Value = "time_icu29"
$ahaGlucose.Range("A141").Value = "BG...140"
$ahaGlucose.Range("B141").Value = "glucose_icu29"
$ahaGlucose.Range("A142").Value = "Time 30:"
$ahaGlucose.Range("B142").Value = "time_icu30"
$ahaGlucose.Range("A143").Value = "BG...142"
$ahaGlucose.Range("B143").Value = "glucose_icu30"
$ahaGlucose.Range("A144").Value = "Average BG"
$ahaGlucose.Range("A145").Value = "HyperglycemiaNumber of BF >140mg/dL"
$ahaGlucose.Range("A146").Value = "Hyperglycemia BG > 180mg/dL:"

# --- Populate "aha patient" ---
$ahaPatient.Range("A1").Value = "variable"
$ahaPatient.Range("B1").Value = "new_var"
$ahaPatient.Range("A2").Value = "Record ID   Group: [screeningrandomiza_arm_1][bllinded_group]"
$ahaPatient.Range("B2").Value = "record_id"
$ahaPatient.Range("A3").Value = "...2"
$ahaPatient.Range("A4").Value = "HbA1c"
$ahaPatient.Range("B4").Value = "hba1c"
$ahaPatient.Range("A5").Value = "Length of surgery"
$ahaPatient.Range("B5").Value = "duration_surgery"
$ahaPatient.Range("A6").Value = "# of vessels"
$ahaPatient.Range("B6").Value = "n_vessels"
$ahaPatient.Range("A7").Value = "APACHE II"
$ahaPatient.Range("B7").Value = "apache_ii"
$ahaPatient.Range("A8").Value = "Race"
$ahaPatient.Range("B8").Value = "race"
$ahaPatient.Range("A9").Value = "Sex"
$ahaPatient.Range("B9").Value = "sex"
$ahaPatient.Range("A10").Value = "BMI"
$ahaPatient.Range("B10").Value = "bmi"
$ahaPatient.Range("A11").Value = "EF"
$ahaPatient.Range("B11").Value = "ef"
$ahaPatient.Range("A12").Value = "pressor?"
$ahaPatient.Range("B12").Value = "pressor"
$ahaPatient.Range("A13").Value = ">=2 pressors"
$ahaPatient.Range("B13").Value = "ge2_pressor"
$ahaPatient.Range("A14").Value = "AF"
$ahaPatient.Range("B14").Value = "afib"
$ahaPatient.Range("A15").Value = "Age"
$ahaPatient.Range("B15").Value = "age"
$ahaPatient.Range("A16").Value = "Any cardiac complication? (MI, Arrhythmia, Cardiac arrest, Acute heart failure, vasopressors, inotropes)"
$ahaPatient.Range("B16").Value = "cardiac_complication"
$ahaPatient.Range("A17").Value = "MI"
$ahaPatient.Range("B17").Value = "mi"
$ahaPatient.Range("A18").Value = "Cardiac arrhythmia"
$ahaPatient.Range("B18").Value = "cardiac_arrhythmia"
$ahaPatient.Range("A19").Value = "Arrhythmias (choice=Frequent PVCs/PACs)"
$ahaPatient.Range("B19").Value = "arr_pvc"
$ahaPatient.Range("A20").Value = "Arrhythmias (choice=VFib/Vtach)"
$ahaPatient.Range("B20").Value = "arr_vtach"
$ahaPatient.Range("A21").Value = "Arrhythmias (choice=Cardiac arrest (Vfib/Vtach))"
$ahaPatient.Range("B21").Value = "arr_arrest_vfib"
$ahaPatient.Range("A22").Value = "Arrhythmias (choice=Cardiac arrest (asystole, PEA))"
$ahaPatient.Range("B22").Value = "arr_arrest_asystolepea"
$ahaPatient.Range("A23").Value = "Arrhythmias (choice=Other)"
$ahaPatient.Range("B23").Value = "arr_other"
$ahaPatient.Range("A24").Value = "Specify (1):"
$ahaPatient.Range("A25").Value = "Cardiac arrest"
$ahaPatient.Range("B25").Value = "cardiacarrest"
$ahaPatient.Range("A26").Value = "Pulmonary edema"
$ahaPatient.Range("B26").Value = "pulmonaryedema"
$ahaPatient.Range("A27").Value = "Heart Failure"
$ahaPatient.Range("B27").Value = "heartfailure"
$ahaPatient.Range("A28").Value = "HTN"
$ahaPatient.Range("B28").Value = "hypertension"
$ahaPatient.Range("A29").Value = "Hyperlipidemia"
$ahaPatient.Range("B29").Value = "hyperlipidemia"
$ahaPatient.Range("A30").Value = "Alcohol"
$ahaPatient.Range("B30").Value = "alcohol"
$ahaPatient.Range("A31").Value = "Smoking?"
$ahaPatient.Range("B31").Value = "smoking"
$ahaPatient.Range("A32").Value = "Never smoked"
$ahaPatient.Range("B32").Value = "neversmoked"

# --- Header formatting (bold, matches existing variable-list sheets) ---
$ahaGlucose.Range("A1:B1").Font.Bold = $true
$ahaPatient.Range("A1:B1").Font.Bold = $true

# --- Column widths (autofit to content, matches existing sheets behavior) ---
$ahaGlucose.Columns.Item(1).AutoFit()
$ahaGlucose.Columns.Item(2).AutoFit()
$ahaPatient.Columns.Item(1).AutoFit()

# --- View/selection state for existing sheets ---
$orToIcu = $wb.Worksheets.Item("or_to_icu")
$orToIcu.Select()
$excel.ActiveWindow.Zoom = 85
$orToIcu.Range("B5").Select()

$dtSurgery = $wb.Worksheets.Item("dt_surgery")
$dtSurgery.Select()
$dtSurgery.Range("B3:B4").Select()

$icu48h = $wb.Worksheets.Item("icu48h")
$icu48h.Select()
$icu48h.Range("B11:B42").Select()

# --- Final selections on the new sheets ---
$ahaGlucose.Range("B10").Select()
$ahaPatient.Range("A14").Select()

# --- Make "aha glucose" the active sheet/tab, matching target workbook view ---
$ahaGlucose.Activate()
$ahaGlucose.Range("B10").Select()
$excel.Windows.Item(1).ScrollWorkbookTabs(1, 1)
